# Katsuno_2010.xlsx — "Updated soil type data"
#
# 1. Add a new "pro_usda_soil_order" column to the "profile" sheet
#    (inserted right before the existing "pro_soil_taxon" column, i.e. at
#    column N) and record the USDA soil order ("Andisols") for both
#    profile data rows.
# 2. Wrap the long bibliographical_reference text on the "metadata" sheet
#    so it displays fully.

$wb = $excel.ActiveWorkbook

# --- profile sheet -------------------------------------------------------
$wsProfile = $wb.Worksheets.Item("profile")

# Insert a new blank column at N; everything from pro_soil_taxon onward
# shifts one column to the right.
$wsProfile.Columns("N").Insert()

$wsProfile.Range("N1").Value = "pro_usda_soil_order"
$wsProfile.Range("N4").Value = "Andisols"
$wsProfile.Range("N5").Value = "Andisols"

$wsProfile.Range("O11").Select() | Out-Null

# --- metadata sheet --------------------------------------------------------
$wsMetadata = $wb.Worksheets.Item("metadata")

# The bibliographical_reference cell holds a long citation string; wrap it.
$wsMetadata.Range("M4").WrapText = $true
$wsMetadata.Rows.Item(4).RowHeight = 409.6

$wsMetadata.Activate() | Out-Null
$wsMetadata.Range("A4").Select() | Out-Null
